$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 192, pushing existing rows 192-256 down to 193-257.
$ws.Rows.Item(192).Insert()

# Copy formatting/content of the (now shifted) old row 192 -- which is now row 193 -- into new row 192,
# then overwrite the changed columns with the new week's data.
$srcRow = $ws.Range("A193:R193")
$dstRow = $ws.Range("A192:R192")
$srcRow.Copy()
$dstRow.PasteSpecial(-4104)  # xlPasteAll

# Update the changed values for the new row (192)
$ws.Cells.Item(192, 4).Value = 44988     # D: Fecha
$ws.Cells.Item(192, 10).Value = 180      # J: Volumen
$ws.Cells.Item(192, 11).Value = 15000    # K: Precio minimo
$ws.Cells.Item(192, 12).Value = 16000    # L: Precio maximo
$ws.Cells.Item(192, 13).Value = 15444    # M: Precio promedio ponderado
$ws.Cells.Item(192, 16).Value = 1544     # P: Precio $/Kg

$ws.Range("A1").Select()
